$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append below the existing data (rows 2..11 already populated).
# Columns: A, B, C, D, E, F, G, H
$rows = @(
    @(8,    8,    1,    6,    13.88, "1000P", 0.01, 162.0),
    @(8,    8,    1,    6,    13.7,  "1000P", 0.01, 198.0),
    @(8,    8,    1,    6,    13.89, "1000P", 0.01, 216.0),
    @(8,    8,    1,    6,    13.27, "1000P", 0.01, 234.0),
    @(8,    8,    1,    6,    14.91, "1000P", 0.01, 252.0),
    @(8,    8,    1,    6,    15.41, "1000P", 0.01, 270.0),
    @(8,    8,    1,    6,    14.43, "1000P", 0.01, 288.0),
    @(8,    8,    1,    6,    13.78, "1000P", 0.01, 306.0),
    @(8,    8,    1,    6,    13.97, "1000P", 0.01, 324.0),
    @(8,    8,    1,    6,    14.4,  "10P",   0.01, 0),
    @(8,    8,    1,    6,    15.55, "10P",   0.01, 72.0),
    @(8,    8,    1,    6,    16.79, "10P",   0.01, 144.0),
    @(8,    8,    1,    6,    18.13, "10P",   0.01, 216.0),
    @(8,    8,    1,    6,    14.32, "10P",   0.01, 288.0),
    @(1.5,  8.0,  1.0,  4.02, 15.65, "1000P", 0.01, 0)
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
